$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing stimuli colours
$ws.Range("A2").Value = "images/YellowIcon.png"
$ws.Range("A3").Value = "images/YellowIcon.png"
$ws.Range("A4").Value = "images/OrangeIcon.png"
$ws.Range("A5").Value = "images/OrangeIcon.png"

# Add two new rows using the original Stop icon
$ws.Range("A6").Value = "images/StopIcon.png"
$ws.Range("B6").Value = 0.5
$ws.Range("A7").Value = "images/StopIcon.png"
$ws.Range("B7").Value = -0.5

# Update selection to match final workbook state
$ws.Range("D9").Select()
